$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1649.9286
$ws.Range("I15").Value = 1649.9286
$ws.Range("K15").Value = 4949.7858
$ws.Range("M15").Value = -4780.7858
$ws.Range("H40").Value = 5093.2
$ws.Range("I40").Value = 1999.8572
$ws.Range("K40").Value = 1999.8572
$ws.Range("M40").Value = -1824.8572
$ws.Range("H42").Value = 76.77778000000001
$ws.Range("I42").Value = 91.40000000000001
$ws.Range("J42").Value = 58.5
$ws.Range("K42").Value = 274.2
$ws.Range("L42").Value = 175.5
$ws.Range("M42").Value = -44.20000000000005
$ws.Range("N42").Value = -635.5
$ws.Range("H43").Value = 100000000
$ws.Range("I43").Value = 100000000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 100000000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -99999931
$ws.Range("N43").ClearContents()
$ws.Range("H53").Value = 75.23077000000001
$ws.Range("I53").Value = 75.75
$ws.Range("K53").Value = 75.75
$ws.Range("M53").Value = 561.25
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H74").Value = 2500
$ws.Range("I74").Value = 2500
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2500
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1564
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2500
$ws.Range("I77").Value = 2500
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 12500
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -7820
$ws.Range("N77").ClearContents()
$ws.Range("H116").Value = 2000
$ws.Range("J116").Value = 2000
$ws.Range("L116").Value = 2000
$ws.Range("N116").Value = -8884

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4023.9473
$ws.Range("I32").Value = 4023.9473
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4023.9473
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3736.9473
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 1000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H122").Value = 3429.7
$ws.Range("I122").Value = 4466
$ws.Range("J122").Value = 2985.5715
$ws.Range("K122").Value = 13398
$ws.Range("L122").Value = 8956.7145
$ws.Range("M122").Value = -10948
$ws.Range("N122").Value = -13856.7145
$ws.Range("H136").Value = 1000
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3090.6667
$ws.Range("I94").Value = 2708.8
$ws.Range("K94").Value = 2708.8
$ws.Range("M94").Value = -2257.8
$ws.Range("H99").Value = 4899.6665
$ws.Range("I99").Value = 4899.6665
$ws.Range("K99").Value = 4899.6665
$ws.Range("M99").Value = -3401.6665
$ws.Range("H107").Value = 10416.846
$ws.Range("I107").Value = 4456.4546
$ws.Range("J107").Value = 43199
$ws.Range("K107").Value = 4456.4546
$ws.Range("L107").Value = 43199
$ws.Range("M107").Value = -2536.4546
$ws.Range("N107").Value = -47039

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 495.66666
$ws.Range("J107").Value = 495.66666
$ws.Range("L107").Value = 495.66666
$ws.Range("N107").Value = -4335.66666

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 9
$ws.Range("K2").Value = 90
$ws.Range("L2").Value = 54
$ws.Range("M2").Value = 23
$ws.Range("N2").Value = -280
$ws.Range("H4").Value = 200619.8
$ws.Range("I4").Value = 775
$ws.Range("K4").Value = 2325
$ws.Range("M4").Value = -2213
$ws.Range("H12").Value = 15.285714
$ws.Range("I12").Value = 8.5
$ws.Range("K12").Value = 25.5
$ws.Range("M12").Value = 147.5
$ws.Range("H68").Value = 600
$ws.Range("I68").Value = 600
$ws.Range("K68").Value = 1800
$ws.Range("M68").Value = -989
$ws.Range("H71").Value = 600
$ws.Range("I71").Value = 600
$ws.Range("K71").Value = 5400
$ws.Range("M71").Value = -1344
$ws.Range("H87").Value = 400
$ws.Range("I87").Value = 400
$ws.Range("K87").Value = 1200
$ws.Range("M87").Value = 48
$ws.Range("H90").Value = 400
$ws.Range("I90").Value = 400
$ws.Range("K90").Value = 3600
$ws.Range("M90").Value = 2640

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 833.3333
$ws.Range("K22").Value = 833.3333
$ws.Range("M22").Value = -538.3333
$ws.Range("H27").Value = 900
$ws.Range("I27").Value = 833.3333
$ws.Range("K27").Value = 833.3333
$ws.Range("M27").Value = -726.3333
$ws.Range("H55").Value = 1387.091
$ws.Range("I55").Value = 677.3333
$ws.Range("J55").Value = 2238.8
$ws.Range("K55").Value = 677.3333
$ws.Range("L55").Value = 2238.8
$ws.Range("M55").Value = -504.3333
$ws.Range("N55").Value = -2584.8

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 35000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 35000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 35000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -36040
$ws.Range("H62").Value = 125
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 125
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 125
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -1373
$ws.Range("H65").Value = 125
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 125
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 625
$ws.Range("M65").Value = 625
$ws.Range("N65").Value = -6865
$ws.Range("H107").Value = 1151
$ws.Range("J107").Value = 3000
$ws.Range("L107").Value = 9000
$ws.Range("N107").Value = -12840
